# Restored from revision #9dc7411cce823db9559364a7ef17392df71c9470.TEST
# Author: admin. Type: SAVE.
#
# Rules sheet, row 10 (rule "R30"): the "From" value in column C
# changes from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
